# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column G (direccion-provincial-nombre) is re-tagged from a
# dimension/reference-area column to a measure column, matching the
# newly curated dimensions:
#   G2: sdmx-dimension:refArea              -> iaest-measure:direccion-provincial-nombre
#   G3: dim                                 -> medida
#   G4: URI-Provincia                       -> xsd:int

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"
